# Add a hypothetical "exceptional tree" bonus-caliper column to the
# View Ave tree-replacement worksheet.
#
# 1. Insert a new (blank) column at D - this shifts the old D:N data
#    block one column to the right (new E:O), matching the old <col>
#    width definitions which move from col 4-8/13 to col 5-9/14.
# 2. Populate the freshly-vacated column C with a header label and, for
#    every data row, a bonus-caliper formula that applies a 1.5x
#    "exceptional tree" multiplier once a tree's caliper (col A) is
#    30" or more.
# 3. Add the literal helper value in B3 and the two new grand-total
#    formulas in row 42.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Shift D:N -> E:O by inserting a new blank column at D.
$ws.Range("D1").EntireColumn.Insert()

# 2. New column header + per-row bonus-caliper formulas in column C.
$ws.Range("C1").Value = "With exceptional bonus"

$ws.Range("C3").Formula = "=IF(A3>=30, A3*1.5, A3)"
$ws.Range("C4:C41").Formula = "=IF(A4>=30, A4*1.5, A4)"

# The category-header rows (text in column A, not a caliper number)
# don't get a bonus-caliper value - clear them back out individually.
$ws.Range("C5").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("C25").ClearContents()
$ws.Range("C31").ClearContents()
$ws.Range("C39").ClearContents()

# 3. B3 literal helper + new grand-total row 42.
$ws.Range("B3").Value = 34
$ws.Range("B42").Formula = "=SUM(B3:B41)"
$ws.Range("C42").Formula = "=SUM(C3:C41)"

# Cosmetic: matches the author's final selection/scroll position.
$ws.Range("C5").Select()
